$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.62"
$ws.Range("E2").Value = "'1.08%"
$ws.Range("D3").Value = "'31.71"
$ws.Range("E3").Value = "'1.36%"
$ws.Range("D4").Value = "'5.104"
$ws.Range("E4").Value = "'-0.28%"
$ws.Range("D5").Value = "'0.07826"
$ws.Range("E5").Value = "'-2.63%"
$ws.Range("D6").Value = "'2.285"
$ws.Range("E6").Value = "'-10.98%"
$ws.Range("D7").Value = "'7.780"
$ws.Range("E7").Value = "'-0.75%"
$ws.Range("E8").Value = "'-0.26%"
$ws.Range("D9").Value = "'0.9187"
$ws.Range("E9").Value = "'0.08%"
$ws.Range("D10").Value = "'0.1766"
$ws.Range("E10").Value = "'1.92%"
$ws.Range("D11").Value = "'0.07503"
$ws.Range("E11").Value = "'2.56%"
$ws.Range("D12").Value = "'0.08977"
$ws.Range("E12").Value = "'7.53%"
$ws.Range("D13").Value = "'0.03044"
$ws.Range("E13").Value = "'0.34%"
$ws.Range("E14").Value = "'0.67%"
$ws.Range("D15").Value = "'0.001520"
$ws.Range("E15").Value = "'1.99%"
$ws.Range("D16").Value = "'0.005838"
$ws.Range("E16").Value = "'-1.24%"
$ws.Range("D17").Value = "'3.463"
$ws.Range("E17").Value = "'-1.18%"
$ws.Range("E18").Value = "'-0.03%"
$ws.Range("E19").Value = "'0.25%"
$ws.Range("E20").Value = "'-0.03%"
$ws.Range("D21").Value = "'4.214"
$ws.Range("E21").Value = "'-8.99%"
$ws.Range("D23").Value = "'0.04593"
$ws.Range("E23").Value = "'0.56%"
$ws.Range("E24").Value = "'-0.78%"
$ws.Range("D25").Value = "'0.004468"
$ws.Range("E25").Value = "'0.42%"
$ws.Range("E26").Value = "'5.74%"
$ws.Range("E27").Value = "'-1.38%"
$ws.Range("D39").Value = "'0.01771"
$ws.Range("E39").Value = "'-3.69%"
$ws.Range("D40").Value = "'0.04791"
$ws.Range("E40").Value = "'6.03%"
$ws.Range("D41").Value = "'0.007357"
$ws.Range("E41").Value = "'4.81%"
$ws.Range("D42").Value = "'0.1360"
$ws.Range("E42").Value = "'1.28%"
$ws.Range("E43").Value = "'-2.41%"
$ws.Range("D44").Value = "'0.01024"
$ws.Range("E44").Value = "'4.35%"
$ws.Range("D45").Value = "'0.00006295"
$ws.Range("E45").Value = "'-2.99%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.11%"
$ws.Range("D48").Value = "'0.7325"
$ws.Range("E48").Value = "'-10.73%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.11%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.11%"
